$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = "Anirudh"
$ws.Cells.Item(2, 3).Value = "O"
$ws.Cells.Item(2, 4).Value = "UNSOLD"
$ws.Cells.Item(2, 5).Value = 0

# Row 3
$ws.Cells.Item(3, 1).Value = 45
$ws.Cells.Item(3, 2).Value = "Shubham Panchal "
$ws.Cells.Item(3, 3).Value = "AR"
$ws.Cells.Item(3, 4).Value = "DC"
$ws.Cells.Item(3, 5).Value = 100

# Row 4
$ws.Cells.Item(4, 1).Value = 16
$ws.Cells.Item(4, 2).Value = "Dhananjay (Bachhu)"
$ws.Cells.Item(4, 3).Value = "B"
$ws.Cells.Item(4, 4).Value = "DC"
$ws.Cells.Item(4, 5).Value = 100

# Row 5
$ws.Cells.Item(5, 1).Value = 35
$ws.Cells.Item(5, 2).Value = "Vishal Sawant "
$ws.Cells.Item(5, 3).Value = "'="
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "MI"
$ws.Cells.Item(5, 5).Value = 20

# Row 6
$ws.Cells.Item(6, 1).Value = 12
$ws.Cells.Item(6, 2).Value = "Viraj Ambre"
$ws.Cells.Item(6, 3).Value = "'"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "UNSOLD"
$ws.Cells.Item(6, 5).Value = 0

# Row 7
$ws.Cells.Item(7, 1).Value = 34
$ws.Cells.Item(7, 2).Value = "Ganesh Pandian "
$ws.Cells.Item(7, 3).Value = "O"
$ws.Cells.Item(7, 4).Value = "DC"
$ws.Cells.Item(7, 5).Value = 100

# Row 8
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "Monish Rathod"
$ws.Cells.Item(8, 3).Value = "B"
$ws.Cells.Item(8, 4).Value = "MI"
$ws.Cells.Item(8, 5).Value = 100

# Row 9
$ws.Cells.Item(9, 1).Value = 57
$ws.Cells.Item(9, 2).Value = " Chetan Shrivastav"
$ws.Cells.Item(9, 3).Value = "AR"
$ws.Cells.Item(9, 4).Value = "DC"
$ws.Cells.Item(9, 5).Value = 100
